# Apply updated cryptocurrency price/volume data to Sheet1
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$cell = $ws.Range("D2")
$cell.NumberFormat = "@"
$cell.Value = "26.124.86"
$cell.Style = "Normal"
$cell = $ws.Range("D3")
$cell.NumberFormat = "@"
$cell.Value = "1.652.51"
$cell.Style = "Normal"
$ws.Range("E3").Value = "  -0.58%  "
$ws.Range("E4").Value = "  -0.32%  "
$cell = $ws.Range("D5")
$cell.NumberFormat = "@"
$cell.Value = "218.71"
$cell.Style = "Normal"
$ws.Range("E5").Value = "  +0.00%  "
$cell = $ws.Range("D6")
$cell.NumberFormat = "@"
$cell.Value = "0.5287"
$cell.Style = "Normal"
$ws.Range("E6").Value = "  +1.08%  "
$ws.Range("E7").Value = "  -0.23%  "
$cell = $ws.Range("D8")
$cell.NumberFormat = "@"
$cell.Value = "0.2608"
$cell.Style = "Normal"
$ws.Range("E8").Value = "  -2.30%  "
$cell = $ws.Range("D9")
$cell.NumberFormat = "@"
$cell.Value = "0.06322"
$cell.Style = "Normal"
$ws.Range("E9").Value = "  +0.04%  "
$cell = $ws.Range("D10")
$cell.NumberFormat = "@"
$cell.Value = "20.38"
$cell.Style = "Normal"
$ws.Range("E10").Value = "  -3.25%  "
$cell = $ws.Range("D11")
$cell.NumberFormat = "@"
$cell.Value = "0.07740"
$cell.Style = "Normal"
$cell = $ws.Range("D12")
$cell.NumberFormat = "@"
$cell.Value = "4.484"
$cell.Style = "Normal"
$ws.Range("E12").Value = "  +1.24%  "
$cell = $ws.Range("D13")
$cell.NumberFormat = "@"
$cell.Value = "1.652.94"
$cell.Style = "Normal"
$ws.Range("E13").Value = "  -0.59%  "
$cell = $ws.Range("D14")
$cell.NumberFormat = "@"
$cell.Value = "0.5460"
$cell.Style = "Normal"
$ws.Range("E14").Value = "  -0.24%  "
$cell = $ws.Range("D15")
$cell.NumberFormat = "@"
$cell.Value = "0.0₅8131"
$cell.Style = "Normal"
$ws.Range("E15").Value = "  -0.87%  "
$cell = $ws.Range("D16")
$cell.NumberFormat = "@"
$cell.Value = "65.21"
$cell.Style = "Normal"
$cell = $ws.Range("D17")
$cell.NumberFormat = "@"
$cell.Value = "26.129.24"
$cell.Style = "Normal"
$ws.Range("E17").Value = "  -0.45%  "
$cell = $ws.Range("D18")
$cell.NumberFormat = "@"
$cell.Value = "1.002"
$cell.Style = "Normal"
$ws.Range("E18").Value = "  -0.35%  "
$cell = $ws.Range("D19")
$cell.NumberFormat = "@"
$cell.Value = "4.543"
$cell.Style = "Normal"
$ws.Range("E19").Value = "  -2.58%  "
$ws.Range("E20").Value = "  +0.48%  "
$ws.Range("E21").Value = "  -1.24%  "
$cell = $ws.Range("D22")
$cell.NumberFormat = "@"
$cell.Value = "5.984"
$cell.Style = "Normal"
$ws.Range("E22").Value = "  -1.66%  "
$ws.Range("E23").Value = "  -0.41%  "
$cell = $ws.Range("D24")
$cell.NumberFormat = "@"
$cell.Value = "140.14"
$cell.Style = "Normal"
$ws.Range("E24").Value = "  +0.91%  "
$cell = $ws.Range("D25")
$cell.NumberFormat = "@"
$cell.Value = "0.1239"
$cell.Style = "Normal"
$ws.Range("E25").Value = "  +0.14%  "
$cell = $ws.Range("D26")
$cell.NumberFormat = "@"
$cell.Value = "7.264"
$cell.Style = "Normal"
$cell = $ws.Range("D27")
$cell.NumberFormat = "@"
$cell.Value = "16.15"
$cell.Style = "Normal"
$ws.Range("E27").Value = "  -0.01%  "
$cell = $ws.Range("D28")
$cell.NumberFormat = "@"
$cell.Value = "1.436"
$cell.Style = "Normal"
$ws.Range("E28").Value = "  +1.40%  "
$cell = $ws.Range("D29")
$cell.NumberFormat = "@"
$cell.Value = "0.05943"
$cell.Style = "Normal"
$ws.Range("E29").Value = "  -1.01%  "
$ws.Range("E30").Value = "  -0.36%  "
$cell = $ws.Range("D31")
$cell.NumberFormat = "@"
$cell.Value = "3.504"
$cell.Style = "Normal"
$cell = $ws.Range("D32")
$cell.NumberFormat = "@"
$cell.Value = "3.231"
$cell.Style = "Normal"
$ws.Range("E32").Value = "  -2.37%  "
$cell = $ws.Range("D33")
$cell.NumberFormat = "@"
$cell.Value = "1.547"
$cell.Style = "Normal"
$ws.Range("E33").Value = "  -5.31%  "
$cell = $ws.Range("D34")
$cell.NumberFormat = "@"
$cell.Value = "2.411"
$cell.Style = "Normal"
$ws.Range("E34").Value = "  -0.06%  "
$cell = $ws.Range("D35")
$cell.NumberFormat = "@"
$cell.Value = "0.9440"
$cell.Style = "Normal"
$ws.Range("E35").Value = "  -3.67%  "
$cell = $ws.Range("D36")
$cell.NumberFormat = "@"
$cell.Value = "2.757"
$cell.Style = "Normal"
$ws.Range("E36").Value = "  -1.04%  "
$cell = $ws.Range("D37")
$cell.NumberFormat = "@"
$cell.Value = "0.5634"
$cell.Style = "Normal"
$ws.Range("E37").Value = "  -4.06%  "
$cell = $ws.Range("D38")
$cell.NumberFormat = "@"
$cell.Value = "0.01606"
$cell.Style = "Normal"
$ws.Range("E38").Value = "  +1.15%  "
$cell = $ws.Range("D39")
$cell.NumberFormat = "@"
$cell.Value = "5.853"
$cell.Style = "Normal"
$ws.Range("E39").Value = "  -1.53%  "
$cell = $ws.Range("D40")
$cell.NumberFormat = "@"
$cell.Value = "0.8461"
$cell.Style = "Normal"
$ws.Range("E40").Value = "  -1.81%  "
$ws.Range("E41").Value = "  -0.19%  "
$cell = $ws.Range("D42")
$cell.NumberFormat = "@"
$cell.Value = "100.87"
$cell.Style = "Normal"
$ws.Range("E42").Value = "  +1.25%  "
$cell = $ws.Range("D43")
$cell.NumberFormat = "@"
$cell.Value = "1.008.37"
$cell.Style = "Normal"
$ws.Range("E43").Value = "  -2.40%  "
$cell = $ws.Range("D44")
$cell.NumberFormat = "@"
$cell.Value = "1.798.71"
$cell.Style = "Normal"
$ws.Range("E44").Value = "  -0.25%  "
$cell = $ws.Range("D45")
$cell.NumberFormat = "@"
$cell.Value = "56.78"
$cell.Style = "Normal"
$ws.Range("E45").Value = "  -0.63%  "
$cell = $ws.Range("D46")
$cell.NumberFormat = "@"
$cell.Value = "0.0₈106"
$cell.Style = "Normal"
$ws.Range("E46").Value = "  -2.76%  "
$ws.Range("E47").Value = "  +0.00%  "
$cell = $ws.Range("D48")
$cell.NumberFormat = "@"
$cell.Value = "0.4289"
$cell.Style = "Normal"
$ws.Range("E48").Value = "  +1.37%  "
$ws.Range("E49").Value = "  -0.63%  "
$cell = $ws.Range("D50")
$cell.NumberFormat = "@"
$cell.Value = "1.469"
$cell.Style = "Normal"
$ws.Range("E50").Value = "  +0.42%  "
$cell = $ws.Range("D51")
$cell.NumberFormat = "@"
$cell.Value = "7.737"
$cell.Style = "Normal"
$ws.Range("E51").Value = "  -4.58%  "
